$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
# Tabelle1 is the active/selected sheet in this workbook.

# Add new row 13 data, reusing the formatting from row 2 (A2:C2) so that
# no new cell styles get created in styles.xml
$ws.Range("A2:C2").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = 41755
$ws.Range("B13").Value = "11"
$ws.Range("C13").Value = "JEB"
$ws.Range("D13").Value = "Initial Base Project "
$ws.Range("E13").Value = "Done"

# Update selection on sheet1
$ws.Range("A14").Select()
